$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "51.113.13"
$ws.Range("E2").Value = "  +0.51%  "

# Row 3
$ws.Range("D3").Value = "2.957.70"
$ws.Range("E3").Value = "  +1.18%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.17%  "

# Row 5
$ws.Range("D5").Value = "'379.33"
$ws.Range("E5").Value = "  +1.74%  "

# Row 6
$ws.Range("D6").Value = "'102.12"
$ws.Range("E6").Value = "  +1.37%  "

# Row 7
$ws.Range("E7").Value = "  +2.05%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("D9").Value = "'0.591"
$ws.Range("E9").Value = "  +2.08%  "

# Row 10
$ws.Range("D10").Value = "'36.34"
$ws.Range("E10").Value = "  +1.60%  "

# Row 11
$ws.Range("E11").Value = "  -1.13%  "

# Row 12
$ws.Range("E12").Value = "  +2.29%  "

# Row 13
$ws.Range("D13").Value = "'7.84"
$ws.Range("E13").Value = "  +6.62%  "

# Row 14
$ws.Range("D14").Value = "3.422.82"
$ws.Range("E14").Value = "  +0.61%  "

# Row 15
$ws.Range("D15").Value = "'18.26"
$ws.Range("E15").Value = "  +2.72%  "

# Row 16
$ws.Range("D16").Value = "2.965.06"
$ws.Range("E16").Value = "  +0.55%  "

# Row 17
$ws.Range("D17").Value = "'11.23"
$ws.Range("E17").Value = "  +1.50%  "

# Row 18
$ws.Range("D18").Value = "'0.994"
$ws.Range("E18").Value = "  +2.98%  "

# Row 19
$ws.Range("D19").Value = "51.177.94"
$ws.Range("E19").Value = "  +0.69%  "

# Row 20
$ws.Range("D20").Value = "'3.13"
$ws.Range("E20").Value = "  +0.64%  "

# Row 21
$ws.Range("D21").Value = "'12.35"
$ws.Range("E21").Value = "  -0.58%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0961"
$ws.Range("E22").Value = "  +1.18%  "

# Row 23
$ws.Range("D23").Value = "'70.33"
$ws.Range("E23").Value = "  +2.99%  "

# Row 24
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "'266.89"
$ws.Range("E24").Value = "  +1.66%  "

# Row 25
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "'3.22"
$ws.Range("E25").Value = "  +4.18%  "

# Row 26
$ws.Range("D26").Value = "'7.84"
$ws.Range("E26").Value = "  -1.55%  "

# Row 27
$ws.Range("D27").Value = "'7.24"
$ws.Range("E27").Value = "  -4.20%  "

# Row 28
$ws.Range("E28").Value = "  +0.05%  "

# Row 29
$ws.Range("D29").Value = "'25.84"
$ws.Range("E29").Value = "  +1.55%  "

# Row 30
$ws.Range("D30").Value = "'0.164"
$ws.Range("E30").Value = "  -1.25%  "

# Row 31
$ws.Range("D31").Value = "'0.109"
$ws.Range("E31").Value = "  -0.18%  "

# Row 32
$ws.Range("D32").Value = "'10.28"
$ws.Range("E32").Value = "  +3.98%  "

# Row 33
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").Value = "'51.20"
$ws.Range("E33").Value = "  +1.29%  "

# Row 34
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "'34.33"
$ws.Range("E34").Value = "  +5.33%  "

# Row 35
$ws.Range("D35").Value = "'2.00"
$ws.Range("E35").Value = "  -0.93%  "

# Row 36
$ws.Range("E36").Value = "  -0.67%  "

# Row 37
$ws.Range("E37").Value = "  -0.11%  "

# Row 38
$ws.Range("D38").Value = "'3.23"
$ws.Range("E38").Value = "  +4.64%  "

# Row 39
$ws.Range("D39").Value = "'0.116"
$ws.Range("E39").Value = "  +1.34%  "

# Row 40
$ws.Range("E40").Value = "  +3.43%  "

# Row 41
$ws.Range("D41").Value = "'16.42"
$ws.Range("E41").Value = "  +2.13%  "

# Row 42
$ws.Range("E42").Value = "  +4.11%  "

# Row 43
$ws.Range("D43").Value = "'2.49"
$ws.Range("E43").Value = "  +0.60%  "

# Row 44
$ws.Range("E44").Value = "  +8.62%  "

# Row 45
$ws.Range("D45").Value = "'21.43"
$ws.Range("E45").Value = "  +2.35%  "

# Row 46
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "'0.274"
$ws.Range("E46").Value = "  -0.67%  "

# Row 47
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "'2.02"
$ws.Range("E47").Value = "  +0.38%  "

# Row 48
$ws.Range("D48").Value = "'2.38"
$ws.Range("E48").Value = "  +3.06%  "

# Row 49
$ws.Range("D49").Value = "2.036.01"
$ws.Range("E49").Value = "  +2.58%  "

# Row 50
$ws.Range("E50").Value = "  -3.58%  "

# Row 51
$ws.Range("E51").Value = "  +6.75%  "

